# "minor adjustments to Scenario1"
# - Rows 142-150: convert the formula-driven A/B columns into static values and
#   apply the "wrap text" style (s="1") across the whole row A:J, matching the
#   formatting already used by rows 151-154.
# - Row 154: the "bonus" marker in column G moves off this row (G154: 1 -> 0).
# - Two new rows (155, 156) are appended, continuing the same table pattern,
#   with row 155 carrying the "bonus" marker (G155 = 1) that used to sit on 154.
# - The active selection moves to L155 and the view scrolls down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the table rows that get rewritten as plain values (no formulas).
$rowData = @{
    142 = @(141, 8,    100, 0, 0, 0, 0, 0, 0, 0)
    143 = @(142, 7,    100, 0, 0, 0, 0, 0, 0, 0)
    144 = @(143, 6,    100, 0, 0, 0, 0, 0, 0, 0)
    145 = @(144, 5,    100, 0, 0, 0, 0, 0, 0, 0)
    146 = @(145, 4,    100, 0, 0, 0, 0, 0, 0, 0)
    147 = @(146, 3,    100, 0, 0, 0, 0, 0, 0, 0)
    148 = @(147, 2,    100, 0, 0, 0, 0, 0, 0, 0)
    149 = @(148, 1,    100, 0, 0, 0, 0, 0, 0, 0)
    150 = @(149, 0,   5000, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in 142..150) {
    $values = $rowData[$r]
    $rng = $ws.Range("A$r`:J$r")

    # Apply the wrap-text formatting used by the rest of the "finished" rows.
    $rng.WrapText = $true

    # Build a 1-row x 10-column array so Excel treats this as a single
    # range assignment (this clears any existing formulas in A/B).
    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rng.Value = $arr
}

# Row 154's bonus marker (column G) is cleared; it moves to the new row 155.
$ws.Cells.Item(154, 7).Value = 0

# Append the two new rows, formatted the same way as the rest of the table.
$newRows = @{
    155 = @(154, 0, 100, 0, 0, 0, 1, 0, 0, 0)
    156 = @(155, 0, 100, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in 155..156) {
    $values = $newRows[$r]
    $rng = $ws.Range("A$r`:J$r")

    $arr = New-Object 'object[,]' 1,10
    for ($i = 0; $i -lt 10; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rng.Value = $arr
    $rng.WrapText = $true
}

# Update the view: scroll position and active selection.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 135
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scrolling the window isn't essential to the data change; ignore if unsupported.
}
$ws.Range("L155").Select()
